# Update the arithmetic answer table: replace each "old" equation text
# with its "new" equation text (100 cell replacements total).
$d = $word.ActiveDocument

$d.Content.Find.Execute("41-41=0", $true, $true, $false, $false, $false, $true, 1, $false, "61+37=98", 2) | Out-Null
$d.Content.Find.Execute("30-11=19", $true, $true, $false, $false, $false, $true, 1, $false, "91-55=36", 2) | Out-Null
$d.Content.Find.Execute("94-43=51", $true, $true, $false, $false, $false, $true, 1, $false, "90-81=9", 2) | Out-Null
$d.Content.Find.Execute("23+74=97", $true, $true, $false, $false, $false, $true, 1, $false, "32+54=86", 2) | Out-Null
$d.Content.Find.Execute("54-14=40", $true, $true, $false, $false, $false, $true, 1, $false, "46+21=67", 2) | Out-Null
$d.Content.Find.Execute("71+8=79", $true, $true, $false, $false, $false, $true, 1, $false, "89+8=97", 2) | Out-Null
$d.Content.Find.Execute("25+25=50", $true, $true, $false, $false, $false, $true, 1, $false, "42+52=94", 2) | Out-Null
$d.Content.Find.Execute("22+65=87", $true, $true, $false, $false, $false, $true, 1, $false, "56-48=8", 2) | Out-Null
$d.Content.Find.Execute("73-56=17", $true, $true, $false, $false, $false, $true, 1, $false, "65+18=83", 2) | Out-Null
$d.Content.Find.Execute("18+12=30", $true, $true, $false, $false, $false, $true, 1, $false, "52-10=42", 2) | Out-Null
$d.Content.Find.Execute("43+21=64", $true, $true, $false, $false, $false, $true, 1, $false, "11+18=29", 2) | Out-Null
$d.Content.Find.Execute("70-27=43", $true, $true, $false, $false, $false, $true, 1, $false, "89-70=19", 2) | Out-Null
$d.Content.Find.Execute("8-0=8", $true, $true, $false, $false, $false, $true, 1, $false, "93-17=76", 2) | Out-Null
$d.Content.Find.Execute("44-38=6", $true, $true, $false, $false, $false, $true, 1, $false, "82+3=85", 2) | Out-Null
$d.Content.Find.Execute("41-28=13", $true, $true, $false, $false, $false, $true, 1, $false, "71+17=88", 2) | Out-Null
$d.Content.Find.Execute("68+2=70", $true, $true, $false, $false, $false, $true, 1, $false, "43+33=76", 2) | Out-Null
$d.Content.Find.Execute("48+6=54", $true, $true, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("32+30=62", $true, $true, $false, $false, $false, $true, 1, $false, "29+58=87", 2) | Out-Null
$d.Content.Find.Execute("72-21=51", $true, $true, $false, $false, $false, $true, 1, $false, "25+19=44", 2) | Out-Null
$d.Content.Find.Execute("67+5=72", $true, $true, $false, $false, $false, $true, 1, $false, "51+0=51", 2) | Out-Null
$d.Content.Find.Execute("76-24=52", $true, $true, $false, $false, $false, $true, 1, $false, "77-54=23", 2) | Out-Null
$d.Content.Find.Execute("94-49=45", $true, $true, $false, $false, $false, $true, 1, $false, "89+1=90", 2) | Out-Null
$d.Content.Find.Execute("12+62=74", $true, $true, $false, $false, $false, $true, 1, $false, "89-59=30", 2) | Out-Null
$d.Content.Find.Execute("98-39=59", $true, $true, $false, $false, $false, $true, 1, $false, "15+8=23", 2) | Out-Null
$d.Content.Find.Execute("14-4=10", $true, $true, $false, $false, $false, $true, 1, $false, "12+9=21", 2) | Out-Null
$d.Content.Find.Execute("36+39=75", $true, $true, $false, $false, $false, $true, 1, $false, "46+49=95", 2) | Out-Null
$d.Content.Find.Execute("66+32=98", $true, $true, $false, $false, $false, $true, 1, $false, "40+40=80", 2) | Out-Null
$d.Content.Find.Execute("25+54=79", $true, $true, $false, $false, $false, $true, 1, $false, "22+18=40", 2) | Out-Null
$d.Content.Find.Execute("24-21=3", $true, $true, $false, $false, $false, $true, 1, $false, "64-1=63", 2) | Out-Null
$d.Content.Find.Execute("49-44=5", $true, $true, $false, $false, $false, $true, 1, $false, "86-29=57", 2) | Out-Null
$d.Content.Find.Execute("34+41=75", $true, $true, $false, $false, $false, $true, 1, $false, "15+20=35", 2) | Out-Null
$d.Content.Find.Execute("63+21=84", $true, $true, $false, $false, $false, $true, 1, $false, "58-22=36", 2) | Out-Null
$d.Content.Find.Execute("92-11=81", $true, $true, $false, $false, $false, $true, 1, $false, "66-53=13", 2) | Out-Null
$d.Content.Find.Execute("16+41=57", $true, $true, $false, $false, $false, $true, 1, $false, "45-22=23", 2) | Out-Null
$d.Content.Find.Execute("19+43=62", $true, $true, $false, $false, $false, $true, 1, $false, "25+17=42", 2) | Out-Null
$d.Content.Find.Execute("24+27=51", $true, $true, $false, $false, $false, $true, 1, $false, "70-52=18", 2) | Out-Null
$d.Content.Find.Execute("69-65=4", $true, $true, $false, $false, $false, $true, 1, $false, "90-5=85", 2) | Out-Null
$d.Content.Find.Execute("74+15=89", $true, $true, $false, $false, $false, $true, 1, $false, "11+18=29", 2) | Out-Null
$d.Content.Find.Execute("49+24=73", $true, $true, $false, $false, $false, $true, 1, $false, "9+84=93", 2) | Out-Null
$d.Content.Find.Execute("25+31=56", $true, $true, $false, $false, $false, $true, 1, $false, "12+65=77", 2) | Out-Null
$d.Content.Find.Execute("79-24=55", $true, $true, $false, $false, $false, $true, 1, $false, "66+10=76", 2) | Out-Null
$d.Content.Find.Execute("59+38=97", $true, $true, $false, $false, $false, $true, 1, $false, "59+34=93", 2) | Out-Null
$d.Content.Find.Execute("32+19=51", $true, $true, $false, $false, $false, $true, 1, $false, "56+32=88", 2) | Out-Null
$d.Content.Find.Execute("20+22=42", $true, $true, $false, $false, $false, $true, 1, $false, "75-75=0", 2) | Out-Null
$d.Content.Find.Execute("31+33=64", $true, $true, $false, $false, $false, $true, 1, $false, "32+47=79", 2) | Out-Null
$d.Content.Find.Execute("56-43=13", $true, $true, $false, $false, $false, $true, 1, $false, "16+62=78", 2) | Out-Null
$d.Content.Find.Execute("69-52=17", $true, $true, $false, $false, $false, $true, 1, $false, "52-4=48", 2) | Out-Null
$d.Content.Find.Execute("74-73=1", $true, $true, $false, $false, $false, $true, 1, $false, "40+36=76", 2) | Out-Null
$d.Content.Find.Execute("59+5=64", $true, $true, $false, $false, $false, $true, 1, $false, "14-0=14", 2) | Out-Null
$d.Content.Find.Execute("25+45=70", $true, $true, $false, $false, $false, $true, 1, $false, "85-79=6", 2) | Out-Null
$d.Content.Find.Execute("35+39=74", $true, $true, $false, $false, $false, $true, 1, $false, "10+45=55", 2) | Out-Null
$d.Content.Find.Execute("54+5=59", $true, $true, $false, $false, $false, $true, 1, $false, "78+16=94", 2) | Out-Null
$d.Content.Find.Execute("13+10=23", $true, $true, $false, $false, $false, $true, 1, $false, "54-46=8", 2) | Out-Null
$d.Content.Find.Execute("52-52=0", $true, $true, $false, $false, $false, $true, 1, $false, "99-75=24", 2) | Out-Null
$d.Content.Find.Execute("19-7=12", $true, $true, $false, $false, $false, $true, 1, $false, "42-39=3", 2) | Out-Null
$d.Content.Find.Execute("0+4=4", $true, $true, $false, $false, $false, $true, 1, $false, "76-73=3", 2) | Out-Null
$d.Content.Find.Execute("11+74=85", $true, $true, $false, $false, $false, $true, 1, $false, "0+36=36", 2) | Out-Null
$d.Content.Find.Execute("73+18=91", $true, $true, $false, $false, $false, $true, 1, $false, "98-11=87", 2) | Out-Null
$d.Content.Find.Execute("34+62=96", $true, $true, $false, $false, $false, $true, 1, $false, "60+19=79", 2) | Out-Null
$d.Content.Find.Execute("18+56=74", $true, $true, $false, $false, $false, $true, 1, $false, "97-87=10", 2) | Out-Null
$d.Content.Find.Execute("41+30=71", $true, $true, $false, $false, $false, $true, 1, $false, "68+21=89", 2) | Out-Null
$d.Content.Find.Execute("94-21=73", $true, $true, $false, $false, $false, $true, 1, $false, "26+30=56", 2) | Out-Null
$d.Content.Find.Execute("87+1=88", $true, $true, $false, $false, $false, $true, 1, $false, "77+19=96", 2) | Out-Null
$d.Content.Find.Execute("50-29=21", $true, $true, $false, $false, $false, $true, 1, $false, "1+84=85", 2) | Out-Null
$d.Content.Find.Execute("75-30=45", $true, $true, $false, $false, $false, $true, 1, $false, "3+96=99", 2) | Out-Null
$d.Content.Find.Execute("77+4=81", $true, $true, $false, $false, $false, $true, 1, $false, "19-1=18", 2) | Out-Null
$d.Content.Find.Execute("55+14=69", $true, $true, $false, $false, $false, $true, 1, $false, "20-3=17", 2) | Out-Null
$d.Content.Find.Execute("68-45=23", $true, $true, $false, $false, $false, $true, 1, $false, "71-69=2", 2) | Out-Null
$d.Content.Find.Execute("72-23=49", $true, $true, $false, $false, $false, $true, 1, $false, "89-79=10", 2) | Out-Null
$d.Content.Find.Execute("25+33=58", $true, $true, $false, $false, $false, $true, 1, $false, "8+29=37", 2) | Out-Null
$d.Content.Find.Execute("13+13=26", $true, $true, $false, $false, $false, $true, 1, $false, "53+9=62", 2) | Out-Null
$d.Content.Find.Execute("46-40=6", $true, $true, $false, $false, $false, $true, 1, $false, "98-8=90", 2) | Out-Null
$d.Content.Find.Execute("51+46=97", $true, $true, $false, $false, $false, $true, 1, $false, "89-89=0", 2) | Out-Null
$d.Content.Find.Execute("79+17=96", $true, $true, $false, $false, $false, $true, 1, $false, "87+6=93", 2) | Out-Null
$d.Content.Find.Execute("69-26=43", $true, $true, $false, $false, $false, $true, 1, $false, "33+8=41", 2) | Out-Null
$d.Content.Find.Execute("66+8=74", $true, $true, $false, $false, $false, $true, 1, $false, "84-42=42", 2) | Out-Null
$d.Content.Find.Execute("41-39=2", $true, $true, $false, $false, $false, $true, 1, $false, "89-20=69", 2) | Out-Null
$d.Content.Find.Execute("86-80=6", $true, $true, $false, $false, $false, $true, 1, $false, "49+49=98", 2) | Out-Null
$d.Content.Find.Execute("30+33=63", $true, $true, $false, $false, $false, $true, 1, $false, "50+8=58", 2) | Out-Null
$d.Content.Find.Execute("9-8=1", $true, $true, $false, $false, $false, $true, 1, $false, "54+24=78", 2) | Out-Null
$d.Content.Find.Execute("23+69=92", $true, $true, $false, $false, $false, $true, 1, $false, "94-72=22", 2) | Out-Null
$d.Content.Find.Execute("83-33=50", $true, $true, $false, $false, $false, $true, 1, $false, "29-24=5", 2) | Out-Null
$d.Content.Find.Execute("20+33=53", $true, $true, $false, $false, $false, $true, 1, $false, "60+25=85", 2) | Out-Null
$d.Content.Find.Execute("7+70=77", $true, $true, $false, $false, $false, $true, 1, $false, "41-33=8", 2) | Out-Null
$d.Content.Find.Execute("80-73=7", $true, $true, $false, $false, $false, $true, 1, $false, "71-68=3", 2) | Out-Null
$d.Content.Find.Execute("13+30=43", $true, $true, $false, $false, $false, $true, 1, $false, "56+33=89", 2) | Out-Null
$d.Content.Find.Execute("51+3=54", $true, $true, $false, $false, $false, $true, 1, $false, "57-26=31", 2) | Out-Null
$d.Content.Find.Execute("61+3=64", $true, $true, $false, $false, $false, $true, 1, $false, "12-12=0", 2) | Out-Null
$d.Content.Find.Execute("18+37=55", $true, $true, $false, $false, $false, $true, 1, $false, "97+2=99", 2) | Out-Null
$d.Content.Find.Execute("71-20=51", $true, $true, $false, $false, $false, $true, 1, $false, "92-28=64", 2) | Out-Null
$d.Content.Find.Execute("94-80=14", $true, $true, $false, $false, $false, $true, 1, $false, "42-1=41", 2) | Out-Null
$d.Content.Find.Execute("28+12=40", $true, $true, $false, $false, $false, $true, 1, $false, "48+23=71", 2) | Out-Null
$d.Content.Find.Execute("51-22=29", $true, $true, $false, $false, $false, $true, 1, $false, "20+75=95", 2) | Out-Null
$d.Content.Find.Execute("74-60=14", $true, $true, $false, $false, $false, $true, 1, $false, "34-28=6", 2) | Out-Null
$d.Content.Find.Execute("36+62=98", $true, $true, $false, $false, $false, $true, 1, $false, "40+32=72", 2) | Out-Null
$d.Content.Find.Execute("91-20=71", $true, $true, $false, $false, $false, $true, 1, $false, "58+12=70", 2) | Out-Null
$d.Content.Find.Execute("17-7=10", $true, $true, $false, $false, $false, $true, 1, $false, "67-36=31", 2) | Out-Null
$d.Content.Find.Execute("52+11=63", $true, $true, $false, $false, $false, $true, 1, $false, "7+77=84", 2) | Out-Null
$d.Content.Find.Execute("96-58=38", $true, $true, $false, $false, $false, $true, 1, $false, "11+2=13", 2) | Out-Null
$d.Content.Find.Execute("82-60=22", $true, $true, $false, $false, $false, $true, 1, $false, "69-65=4", 2) | Out-Null
